$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ID", "Имя", "Фамилия"),
    @(1, "Алмас", "Нагимов"),
    @(2, "Айжан", "Досмухамедова"),
    @(3, "Бекжан", "Султанов"),
    @(4, "Гульнар", "Аманжолова"),
    @(5, "Данияр", "Кенжебаев"),
    @(6, "Ержан", "Калыбеков"),
    @(7, "Жанна", "Мухамеджанова"),
    @(8, "Илья", "Кузнецов"),
    @(9, "Камилла", "Сафиуллина"),
    @(10, "Лайла", "Абдуллаева")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Header formatting: bold, centered, wrap text
$header = $ws.Range("A1:C1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108
$header.WrapText = $true

# Body formatting: vertical center, wrap text
$body = $ws.Range("A2:C11")
$body.VerticalAlignment = -4108
$body.WrapText = $true
